$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell F1, copying the header style used by the
# existing header cells (B1:E1 -> style index 1: bold, bordered, centered).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in the time_taken values for each data row (F2:F22), plain/default
# styled text cells matching the rest of the data rows.
$timestamps = @(
    "2021-10-05 10:51:09.803739",
    "2021-10-05 10:51:09.803749",
    "2021-10-05 10:51:09.803753",
    "2021-10-05 10:51:09.803755",
    "2021-10-05 10:51:09.803758",
    "2021-10-05 10:51:09.803761",
    "2021-10-05 10:51:09.803764",
    "2021-10-05 10:51:09.803767",
    "2021-10-05 10:51:09.803769",
    "2021-10-05 10:51:09.803772",
    "2021-10-05 10:51:09.803775",
    "2021-10-05 10:51:09.803778",
    "2021-10-05 10:51:09.803780",
    "2021-10-05 10:51:09.803783",
    "2021-10-05 10:51:09.803785",
    "2021-10-05 10:51:09.803788",
    "2021-10-05 10:51:09.803791",
    "2021-10-05 10:51:09.803794",
    "2021-10-05 10:51:09.803796",
    "2021-10-05 10:51:09.803799",
    "2021-10-05 10:51:09.803801"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
